# Apply updated counts (column F) for matching rows on both the "展览"
# sheet and the "全部类型" sheet (which mirrors the same event rows).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (row, newValue) pairs for column F.
$updates = @{
    "展览" = @(
        @{ Row = 7;  Value = 625 },
        @{ Row = 8;  Value = 7943 },
        @{ Row = 11; Value = 6824 },
        @{ Row = 14; Value = 4858 },
        @{ Row = 16; Value = 5270 },
        @{ Row = 21; Value = 307 },
        @{ Row = 26; Value = 8953 },
        @{ Row = 28; Value = 1595 },
        @{ Row = 30; Value = 38 },
        @{ Row = 32; Value = 818 },
        @{ Row = 33; Value = 71 },
        @{ Row = 37; Value = 1840 },
        @{ Row = 38; Value = 234 },
        @{ Row = 41; Value = 4665 },
        @{ Row = 47; Value = 902 }
    )
    "全部类型" = @(
        @{ Row = 8;  Value = 625 },
        @{ Row = 9;  Value = 7943 },
        @{ Row = 12; Value = 6824 },
        @{ Row = 16; Value = 4858 },
        @{ Row = 18; Value = 5270 },
        @{ Row = 23; Value = 307 },
        @{ Row = 27; Value = 8953 },
        @{ Row = 29; Value = 1595 },
        @{ Row = 30; Value = 38 },
        @{ Row = 32; Value = 818 },
        @{ Row = 33; Value = 71 },
        @{ Row = 37; Value = 1840 },
        @{ Row = 38; Value = 234 },
        @{ Row = 41; Value = 4665 },
        @{ Row = 47; Value = 902 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
